$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.9754265
$ws.Range("H2").Value = 47.950853
$ws.Range("I2").Value = 0.02590932378848207
$ws.Range("J2").Value = 0.01788258088872626
$ws.Range("O2").Value = 0.9919525181111984
$ws.Range("P2").Value = 0.9919525181111983
$ws.Range("Q2").Value = 0.9082530902906666
$ws.Range("R2").Value = 5.449518541743999
$ws.Range("S2").Value = 0.02570081897454316
$ws.Range("T2").Value = 0.01773867114289921
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("G3").Value = 23.9754265
$ws.Range("H3").Value = 47.950853
$ws.Range("I3").Value = 0.02590932378848207
$ws.Range("J3").Value = 0.01788258088872626
$ws.Range("M3").Value = 0.0003073333333333333
$ws.Range("N3").Value = 0.000922
$ws.Range("O3").Value = 0.008047481888801606
$ws.Range("P3").Value = 0.008047481888801606
$ws.Range("Q3").Value = 0.007368447744333333
$ws.Range("R3").Value = 0.04421068646599999
$ws.Range("S3").Value = 0.000208504813938906
$ws.Range("T3").Value = 0.0001439097458270543
$ws.Range("I4").Value = 0.897119731561083
$ws.Range("J4").Value = 0.9287862719314923
$ws.Range("O4").Value = 0.9919525181111984
$ws.Range("P4").Value = 0.9919525181111983
$ws.Range("S4").Value = 0.8899001767692587
$ws.Range("T4").Value = 0.921311881229556
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("I5").Value = 0.897119731561083
$ws.Range("J5").Value = 0.9287862719314923
$ws.Range("M5").Value = 0.0003073333333333333
$ws.Range("N5").Value = 0.000922
$ws.Range("O5").Value = 0.008047481888801606
$ws.Range("P5").Value = 0.008047481888801606
$ws.Range("Q5").Value = 0.2551351751355556
$ws.Range("R5").Value = 2.29621657622
$ws.Range("S5").Value = 0.007219554791824374
$ws.Range("T5").Value = 0.007474390701936248
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1360573333333333
$ws.Range("H6").Value = 0.408172
$ws.Range("I6").Value = 0.0001470319413559032
$ws.Range("J6").Value = 0.000152221876147087
$ws.Range("O6").Value = 0.9919525181111984
$ws.Range("P6").Value = 0.9919525181111983
$ws.Range("Q6").Value = 0.005154214606222222
$ws.Range("R6").Value = 0.046387931456
$ws.Range("S6").Value = 0.0001458487044707662
$ws.Range("T6").Value = 0.0001509968733557139
$ws.Range("D7").Value = "Neutrophils"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1360573333333333
$ws.Range("H7").Value = 0.408172
$ws.Range("I7").Value = 0.0001470319413559032
$ws.Range("J7").Value = 0.000152221876147087
$ws.Range("M7").Value = 0.0003073333333333333
$ws.Range("N7").Value = 0.000922
$ws.Range("O7").Value = 0.008047481888801606
$ws.Range("P7").Value = 0.008047481888801606
$ws.Range("Q7").Value = 0.00004181495377777778
$ws.Range("R7").Value = 0.000376334584
$ws.Range("S7").Value = 0.000001183236885136971
$ws.Range("T7").Value = 0.000001225002791373084
$ws.Range("G8").Value = 70.673643
$ws.Range("H8").Value = 141.347286
$ws.Range("I8").Value = 0.07637429514751654
$ws.Range("J8").Value = 0.05271343713733154
$ws.Range("O8").Value = 0.9919525181111984
$ws.Range("P8").Value = 0.9919525181111983
$ws.Range("Q8").Value = 2.677306059888
$ws.Range("R8").Value = 16.063836359328
$ws.Range("S8").Value = 0.07575967439054691
$ws.Range("T8").Value = 0.05228922670667237
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("G9").Value = 70.673643
$ws.Range("H9").Value = 141.347286
$ws.Range("I9").Value = 0.07637429514751654
$ws.Range("J9").Value = 0.05271343713733154
$ws.Range("M9").Value = 0.0003073333333333333
$ws.Range("N9").Value = 0.000922
$ws.Range("O9").Value = 0.008047481888801606
$ws.Range("P9").Value = 0.008047481888801606
$ws.Range("Q9").Value = 0.021720366282
$ws.Range("R9").Value = 0.130322197692
$ws.Range("S9").Value = 0.0006146207569696277
$ws.Range("T9").Value = 0.0004242104306591575
$ws.Range("G10").Value = 0.289655
$ws.Range("H10").Value = 0.868965
$ws.Range("I10").Value = 0.0003130190481471841
$ws.Range("J10").Value = 0.0003240679973299332
$ws.Range("O10").Value = 0.9919525181111984
$ws.Range("P10").Value = 0.9919525181111983
$ws.Range("Q10").Value = 0.01097290381333333
$ws.Range("R10").Value = 0.09875613431999999
$ws.Range("S10").Value = 0.0003105000330263697
$ws.Range("T10").Value = 0.0003214600659906803
$ws.Range("D11").Value = "Neutrophils"
$ws.Range("G11").Value = 0.289655
$ws.Range("H11").Value = 0.868965
$ws.Range("I11").Value = 0.0003130190481471841
$ws.Range("J11").Value = 0.0003240679973299332
$ws.Range("M11").Value = 0.0003073333333333333
$ws.Range("N11").Value = 0.000922
$ws.Range("O11").Value = 0.008047481888801606
$ws.Range("P11").Value = 0.008047481888801606
$ws.Range("Q11").Value = 0.00008902063666666666
$ws.Range("R11").Value = 0.00080118573
$ws.Range("S11").Value = 0.000002519015120814382
$ws.Range("T11").Value = 0.000002607931339252844
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.1264026666666667
$ws.Range("H12").Value = 0.379208
$ws.Range("I12").Value = 0.0001365985134151518
$ws.Range("J12").Value = 0.0001414201689728462
$ws.Range("O12").Value = 0.9919525181111984
$ws.Range("P12").Value = 0.9919525181111983
$ws.Range("Q12").Value = 0.004788470087111111
$ws.Range("R12").Value = 0.043096230784
$ws.Range("S12").Value = 0.0001354992393524061
$ws.Range("T12").Value = 0.0001402820927243259
$ws.Range("D13").Value = "Neutrophils"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.1264026666666667
$ws.Range("H13").Value = 0.379208
$ws.Range("I13").Value = 0.0001365985134151518
$ws.Range("J13").Value = 0.0001414201689728462
$ws.Range("M13").Value = 0.0003073333333333333
$ws.Range("N13").Value = 0.000922
$ws.Range("O13").Value = 0.008047481888801606
$ws.Range("P13").Value = 0.008047481888801606
$ws.Range("Q13").Value = 0.007368447744333333
$ws.Range("R13").Value = 0.04421068646599999
$ws.Range("S13").Value = 0.000208504813938906
$ws.Range("T13").Value = 0.0001439097458270543
